$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.271.95'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.78%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.880.62'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.70%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.02%  '

$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4676'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.31%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2823'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.82%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06592'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.59'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.40%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07758'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.76%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.74%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.876.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.067'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.84%  '

$ws.Range("E15").Value = '  +0.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '283.22'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.293.02'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.124.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.375'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007240'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.99%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.163'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.335'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.45%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.70'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.35%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.983'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.372'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.85%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09653'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.42%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.366'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.82%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.468'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.106'
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04664'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.35%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7036'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.24%  '

$ws.Range("E36").Value = '  -1.36%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.717'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01863'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.65%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.536'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.80%  '

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.518'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.70%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '71.88'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.00%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8617'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.20%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.952'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.05%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.53%  '

$ws.Range("B46").Value = 'TheSandbox'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4171'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.11%  '

$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '985.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.33%  '

$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.219'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.55%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.168'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.98%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.93%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1145'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.85%  '
